$wb = $excel.ActiveWorkbook

# Sheet ALC, row 5 (item id 5503)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 517.5789
$ws.Range("I5").Value = 525.86664
$ws.Range("J5").Value = 486.5
$ws.Range("K5").Value = 525.86664
$ws.Range("L5").Value = 486.5
$ws.Range("M5").Value = -410.86664
$ws.Range("N5").Value = -716.5

# Sheet ALC, row 28 (item id 27772)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1603.0435
$ws.Range("I28").Value = 672.4211
$ws.Range("J28").Value = 6023.5
$ws.Range("K28").Value = 672.4211
$ws.Range("L28").Value = 6023.5
$ws.Range("M28").Value = -187.4211
$ws.Range("N28").Value = -6993.5

# Sheet ALC, row 63 (item id 10652)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0

# Sheet ALC, row 66 (item id 10652)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0

# Sheet ALC, row 70 (item id 12604)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1838.1177
$ws.Range("I70").Value = 1245.6666
$ws.Range("J70").Value = 3260
$ws.Range("K70").Value = 3736.9998
$ws.Range("L70").Value = 9780
$ws.Range("M70").Value = -3466.9998
$ws.Range("N70").Value = -10320

# Sheet ALC, row 73 (item id 12604)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 1838.1177
$ws.Range("I73").Value = 1245.6666
$ws.Range("J73").Value = 3260
$ws.Range("K73").Value = 3736.9998
$ws.Range("L73").Value = 9780
$ws.Range("M73").Value = -2800.9998
$ws.Range("N73").Value = -11652

# Sheet ALC, row 106 (item id 19903)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 2551.5264
$ws.Range("I106").Value = 2330.625
$ws.Range("K106").Value = 2330.625
$ws.Range("M106").Value = -1699.625

# Sheet ALC, row 132 (item id 44049)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3646.4783
$ws.Range("I132").Value = 3749.8
$ws.Range("K132").Value = 11249.4
$ws.Range("M132").Value = -8719.400000000001

# Sheet ALC, row 136 (item id 42164)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H136").Value = 90999.60000000001
$ws.Range("J136").Value = 90999.60000000001
$ws.Range("L136").Value = 90999.60000000001
$ws.Range("N136").Value = -101199.6

# Sheet ALC, row 137 (item id 44013)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2386.2942
$ws.Range("I137").Value = 2478.6538
$ws.Range("J137").Value = 2086.125
$ws.Range("K137").Value = 7435.9614
$ws.Range("L137").Value = 6258.375
$ws.Range("M137").Value = -4885.9614
$ws.Range("N137").Value = -11358.375

# Sheet ARM, row 74 (item id 44000)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 77011260
$ws.Range("I74").Value = 111236820
$ws.Range("J74").Value = 3763.5
$ws.Range("K74").Value = 111236820
$ws.Range("L74").Value = 3763.5
$ws.Range("M74").Value = -111235946
$ws.Range("N74").Value = -5511.5

# Sheet ARM, row 77 (item id 44000)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 77011260
$ws.Range("I77").Value = 111236820
$ws.Range("J77").Value = 3763.5
$ws.Range("K77").Value = 556184100
$ws.Range("L77").Value = 18817.5
$ws.Range("M77").Value = -556179732
$ws.Range("N77").Value = -27553.5

# Sheet ARM, row 110 (item id 27708)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 17429.174
$ws.Range("I110").Value = 18984.3
$ws.Range("J110").Value = 7061.6665
$ws.Range("K110").Value = 18984.3
$ws.Range("L110").Value = 7061.6665
$ws.Range("M110").Value = -16939.3
$ws.Range("N110").Value = -11151.6665

# Sheet BSM, row 105 (item id 19947)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 26130
$ws.Range("I105").Value = 50910
$ws.Range("K105").Value = 50910
$ws.Range("M105").Value = -49163

# Sheet CRP, row 31 (item id 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 27032024
$ws.Range("I31").Value = 4083.7307
$ws.Range("J31").Value = 90916250
$ws.Range("K31").Value = 4083.7307
$ws.Range("L31").Value = 90916250
$ws.Range("M31").Value = -3788.7307
$ws.Range("N31").Value = -90916840

# Sheet CRP, row 34 (item id 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 27032024
$ws.Range("I34").Value = 4083.7307
$ws.Range("J34").Value = 90916250
$ws.Range("K34").Value = 4083.7307
$ws.Range("L34").Value = 90916250
$ws.Range("M34").Value = -3881.7307
$ws.Range("N34").Value = -90916654

# Sheet CRP, row 92 (item id 18041)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 45000
$ws.Range("J92").Value = 45000
$ws.Range("L92").Value = 45000
$ws.Range("N92").Value = -49992

# Sheet CRP, row 94 (item id 32934)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1685.125
$ws.Range("I94").Value = 1201.4286
$ws.Range("K94").Value = 1201.4286
$ws.Range("M94").Value = -750.4286

# Sheet CRP, row 99 (item id 36198)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 28928.666
$ws.Range("I99").Value = 28928.666
$ws.Range("K99").Value = 28928.666
$ws.Range("M99").Value = -27430.666

# Sheet CRP, row 107 (item id 27689)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1228.6875
$ws.Range("I107").Value = 1344.625
$ws.Range("J107").Value = 1112.75
$ws.Range("K107").Value = 1344.625
$ws.Range("L107").Value = 1112.75
$ws.Range("M107").Value = 575.375
$ws.Range("N107").Value = -4952.75

# Sheet CRP, row 126 (item id 36198)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 28928.666
$ws.Range("I126").Value = 28928.666
$ws.Range("K126").Value = 86785.99800000001
$ws.Range("M126").Value = -84315.99800000001

# Sheet CRP, row 132 (item id 44019)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 5354.9375
$ws.Range("I132").Value = 5399.1665
$ws.Range("K132").Value = 16197.4995
$ws.Range("M132").Value = -13667.4995

# Sheet CUL, row 56 (item id 10146)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 22485.818
$ws.Range("I56").Value = 22485.818
$ws.Range("K56").Value = 22485.818
$ws.Range("M56").Value = -21955.818

# Sheet CUL, row 122 (item id 36078)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1841
$ws.Range("I122").Value = 2324
$ws.Range("K122").Value = 20916
$ws.Range("M122").Value = -18466

# Sheet GSM, row 132 (item id 44008)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4938.125
$ws.Range("I132").Value = 5308.5386
$ws.Range("K132").Value = 15925.6158
$ws.Range("M132").Value = -13395.6158

# Sheet LTW, row 25 (item id 3547)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 41250
$ws.Range("I25").Value = 15000
$ws.Range("J25").Value = 50000
$ws.Range("K25").Value = 15000
$ws.Range("L25").Value = 50000
$ws.Range("M25").Value = -14770
$ws.Range("N25").Value = -50460

# Sheet LTW, row 55 (item id 5284)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 952.3333
$ws.Range("I55").Value = 769.4286
$ws.Range("J55").Value = 1112.375
$ws.Range("K55").Value = 769.4286
$ws.Range("L55").Value = 1112.375
$ws.Range("M55").Value = -596.4286
$ws.Range("N55").Value = -1458.375

# Sheet LTW, row 68 (item id 12563)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2647
$ws.Range("I68").Value = 2499.9375
$ws.Range("J68").Value = 5000
$ws.Range("K68").Value = 2499.9375
$ws.Range("L68").Value = 5000
$ws.Range("M68").Value = -1750.9375
$ws.Range("N68").Value = -6498

# Sheet LTW, row 71 (item id 12563)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2647
$ws.Range("I71").Value = 2499.9375
$ws.Range("J71").Value = 5000
$ws.Range("K71").Value = 12499.6875
$ws.Range("L71").Value = 25000
$ws.Range("M71").Value = -8755.6875
$ws.Range("N71").Value = -32488

# Sheet WVR, row 96 (item id 19977)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 6464.125
$ws.Range("I96").Value = 5294.6665
$ws.Range("J96").Value = 9972.5
$ws.Range("K96").Value = 5294.6665
$ws.Range("L96").Value = 9972.5
$ws.Range("M96").Value = -3921.6665
$ws.Range("N96").Value = -12718.5

# Remove now-empty profit cells per diff (value dropped entirely, not just zero)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M63").ClearContents()
$ws.Range("M66").ClearContents()
